$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("New York -- New York") results were overwritten by a failed API call
# (GitHub API rate limit exceeded), so the previously-fetched data values are
# cleared out and the status message reflects the error.

$ws.Range("B4:H4").ClearContents()

$ws.Range("J4").Value = $false

$ws.Range("K4:L4").ClearContents()

$ws.Range("O4").Value = "An error occurred. ... RateLimitExceededException(403, {'message': ""API rate limit exceeded for 132.145.200.60. (But here's the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)"", 'documentation_url': 'https://developer.github.com/v3/#rate-limiting'})"
